$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (B, C, D, E, G) - F column is left unchanged
$data = @{
    2 = @(1.505614041169197, 1.65323645889881, 157.8057217802531, 0.4998867070740569, 161.4644589873952)
    3 = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 6.48142807727062, 9.016289739122548)
    4 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    5 = @(0.1554434735375247, 1.65323645889881, 16.98373111632243, 6.48142807727062, 25.27383912602938)
    6 = @(0.7287194209349384, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 5.964442013611383)
    7 = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 3.594575437922795)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
